$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.863211154937744
$ws.Range("B1").Value = 3.2465500831604
$ws.Range("C1").Value = 2.870581388473511
$ws.Range("D1").Value = 2.5421302318573
$ws.Range("E1").Value = 1.721330881118774
